$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("C2").Value = 0.7557598745023819
$ws.Range("D2").Value = 0.4578094384047755

# Row 3
$ws.Range("C3").Value = 0.9583265794869448
$ws.Range("D3").Value = 0.3483122001019088

# Row 4
$ws.Range("C4").Value = -0.6894842416093365
$ws.Range("D4").Value = 0.4977272484470396

# Row 5
$ws.Range("C5").Value = -1.695898350113051
$ws.Range("D5").Value = 0.1040125237119105
$ws.Range("G5").Value = "No"

# Row 6
$ws.Range("C6").Value = 0.0140105156756316
$ws.Range("D6").Value = 0.9889478732056229

# Row 7
$ws.Range("C7").Value = -1.407285853708961
$ws.Range("D7").Value = 0.1733146228724503

# Row 8
$ws.Range("C8").Value = -2.016036870015598
$ws.Range("D8").Value = 0.05616913431786719
$ws.Range("G8").Value = "No"

# Row 9
$ws.Range("C9").Value = -1.986347489107073
$ws.Range("D9").Value = 0.05959207241818421

# Row 10
$ws.Range("C10").Value = -1.898559588264946
$ws.Range("D10").Value = 0.07081771948327864
$ws.Range("G10").Value = "No"

# Row 11
$ws.Range("C11").Value = -1.309922261347858
$ws.Range("D11").Value = 0.2037361315814943
